# Generate Report for Handback
# Updates the localization-status workbook:
#  - Marks rows as handed back (status text) for both zh-cn and de-de sheets
#  - Stamps the "Latest Handback DateTime" for each language
#  - Adds "Latest Target File" (F) and "Latest Handback File" (G) hyperlinked
#    entries (mirroring the existing Source File / Handoff File hyperlinks)

$wb = $excel.ActiveWorkbook

# ---- zh-cn sheet -----------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Range("H2").Value = "2016-03-21 10:58:52"
$wsZh.Range("H3").Value = "2016-03-21 10:58:52"

$mdUrlZh = "https://github.com/OpenLocalizationTest/oltest/blob/5188bac50427b906efd533f28625de4b6ba96c7c/e2e/50527308-acf1-477c-ac29-3589133d0d67.md"
$xlfUrlZh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e584c33e8c2f998e3c3789aa24b77cf96e36fe55/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.zh-cn.xlf"

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $mdUrlZh, "", "", "50527308-acf1-477c-ac29-3589133d0d67.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $xlfUrlZh, "", "", "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $mdUrlZh, "", "", "50527308-acf1-477c-ac29-3589133d0d67.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $xlfUrlZh, "", "", "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.zh-cn.xlf")

# ---- de-de sheet -------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Range("H2").Value = "2016-03-21 10:58:58"
$wsDe.Range("H3").Value = "2016-03-21 10:58:58"

$mdUrlDe = "https://github.com/OpenLocalizationTest/oltest/blob/5188bac50427b906efd533f28625de4b6ba96c7c/e2e/50527308-acf1-477c-ac29-3589133d0d67.md"
$xlfUrlDe = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/28fac1fbaa18e9bee999b896c6a1cafbf4b5673b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.de-de.xlf"

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $mdUrlDe, "", "", "50527308-acf1-477c-ac29-3589133d0d67.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $xlfUrlDe, "", "", "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $mdUrlDe, "", "", "50527308-acf1-477c-ac29-3589133d0d67.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $xlfUrlDe, "", "", "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.de-de.xlf")

Write-Host "Handback report generated"
